$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Environment sheet: add a third column ("Name") with the QA engineers'
#    first names, matching the header/data styling already used in A:B.
# ---------------------------------------------------------------------------
$wsEnv = $wb.Worksheets.Item(1)

$wsEnv.Range("C2").Value = "Mahesh"
$wsEnv.Range("C3").Value = "Samba"
$wsEnv.Range("C4").Value = "Kiran"

# Header cell C1 picks up the same style as A1/B1.
$wsEnv.Range("B1").Copy()
$wsEnv.Range("C1").PasteSpecial(-4122)
$wsEnv.Range("C1").Value = "Name"

# Data cells C2:C4 pick up the same style as A2:A4.
$wsEnv.Range("A2").Copy()
$wsEnv.Range("C2:C4").PasteSpecial(-4122)
$wsEnv.Range("C2").Value = "Mahesh"
$wsEnv.Range("C3").Value = "Samba"
$wsEnv.Range("C4").Value = "Kiran"

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "AddEmployee" worksheet right before "LoginPage".
# ---------------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item(2)
$wsAdd = $wb.Worksheets.Add($wsLogin)
$wsAdd.Name = "AddEmployee"

$wsAdd.Range("A1").Value = "First Name"
$wsAdd.Range("B1").Value = "Last Name"

$wsAdd.Range("A2").Value = "Mahesh"
$wsAdd.Range("B2").Value = "Goud"
$wsAdd.Range("A3").Value = "Samba"
$wsAdd.Range("B3").Value = "Uppala"
$wsAdd.Range("A4").Value = "Kiran"
$wsAdd.Range("B4").Value = "Chiramana"

$wsAdd.Range("A5").Value = "Pavan"
$wsAdd.Range("B5").Value = "Nadavala"
$wsAdd.Range("A6").Value = "Keshav"
$wsAdd.Range("B6").Value = "Chenna"

$wsAdd.Range("C1").Value = "UserName"
$wsAdd.Range("D1").Value = "Password"

$wsAdd.Range("C4").Value = "KiranChiramana"
$wsAdd.Range("C2").Value = "MaheshGoud"
$wsAdd.Range("C3").Value = "SambaUppala"
$wsAdd.Range("C5").Value = "PavanNadavala"
$wsAdd.Range("C6").Value = "KeshavChenna"

$wsAdd.Range("D2").Value = "Mahesh@Goud123"
$wsAdd.Range("D3").Value = "Samba@Uppala123"
$wsAdd.Range("D4").Value = "Kiran@Chiramana123"
$wsAdd.Range("D5").Value = "Pavan@Nadavala123"
$wsAdd.Range("D6").Value = "Keshav@Chenna123"

# Bold header row with an accent fill + border, matching the style used on
# the other sheets' header rows.
$wsAdd.Range("A1:D1").Font.Bold = $true
$wsAdd.Range("A1:D1").Interior.ThemeColor = 8
$wsAdd.Range("A1:D1").Borders.LineStyle = 1

# Column widths for the new sheet.
$wsAdd.Columns.Item(1).ColumnWidth = 10
$wsAdd.Columns.Item(2).ColumnWidth = 9.333333333333334
$wsAdd.Columns.Item(3).ColumnWidth = 9.666666666666666
$wsAdd.Columns.Item(4).ColumnWidth = 18.5

$wsAdd.Range("H9").Select()
$wsAdd.Activate()

# ---------------------------------------------------------------------------
# 3. Restore the Environment sheet's selection.
# ---------------------------------------------------------------------------
$wsEnv.Range("E6").Select()
